# İş Takip Güncellemesi - 03.08.2025 19:06:19
# Column I ("İHALELİ/MÜDÜRLÜK") values for all data rows change from
# "Müdürlük" to "İhaleli" (rows 2-122), and E95 (GÖREVLİ PERSONELLER)
# is cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("İş Takip Listesi")

for ($row = 2; $row -le 122; $row++) {
    $ws.Cells.Item($row, 9).Value = "İhaleli"
}

$ws.Cells.Item(95, 5).Value = ""
